# Plant-community diversity sheet: the "m" (morning) and "j" (July?) plot-code
# column headers in row 1 were missing a separator dot between the letter
# prefix and the two-digit plot number (e.g. "m01" -> "m.01", "j11" -> "j.11").
# Fix the 22 header cells B1:W1 accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$headerCols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W")
$headerVals = @("m.01","m.02","m.03","m.04","m.05","m.06","m.07","m.08","m.09","m.10","m.11","m.12","m.13","m.14","m.15","m.16","j.11","j.12","j.13","j.14","j.15","j.16")

for ($i = 0; $i -lt $headerCols.Length; $i++) {
    $ws.Range($headerCols[$i] + "1").Value = $headerVals[$i]
}

# Match the author's final selection/cursor position recorded in the sheet.
$ws.Range("W17").Select()

$wb.Save()
